$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A holds the date label as plain text (matches existing shared-string
# entries like "03-11-2021"). Assigning the literal string via .Value would
# be auto-recognized as a date by Excel, so build it as a text formula result
# and paste the computed value back in place - this keeps the cell a plain
# text/shared-string cell with no numeric/date formatting applied.
$cell = $ws.Cells.Item($newRow, 1)
$cell.Formula = '=""&"04-11-2021"'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = 50000
$ws.Cells.Item($newRow, 3).Value = 175000
$ws.Cells.Item($newRow, 4).Value = 50000
$ws.Cells.Item($newRow, 5).Value = 10000
$ws.Cells.Item($newRow, 6).Value = 40000
$ws.Cells.Item($newRow, 7).Value = 3.18
